# Adds an "address_available_time" column (I) to the package_data sheet,
# corrects a data-entry mistake that had the wrong_address flag and the
# later availability time on the wrong row, and updates which sheet/cell
# is active when the workbook is opened.

$wb = $excel.ActiveWorkbook

$pkg = $wb.Worksheets.Item("package_data")
$dist = $wb.Worksheets.Item("distance_data")

# --- New column header: I1 = "address_available_time" ---------------------
$header = $pkg.Range("I1")
$header.Value2 = "address_available_time"
$header.Font.Bold = $true
$header.Font.Size = 9
$header.Borders.Item(7).LineStyle = 1
$header.Borders.Item(10).LineStyle = 1
$header.HorizontalAlignment = -4108
$header.WrapText = $true

# Header row is taller now that there is more text wrapping in it.
$pkg.Rows.Item(1).RowHeight = 36.5

# --- New column body: I2:I41, all zero except the corrected row ------------
for ($r = 2; $r -le 41; $r++) {
    $cell = $pkg.Range("I" + $r)
    $cell.NumberFormat = "h:mm"
    $cell.Value2 = 0
}

# --- Fix the swapped wrong_address flag / available_time value -------------
# Row 9 (package 8 / address 12) incorrectly had wrong_address = TRUE.
$pkg.Range("H9").Value2 = $false

# Row 10 (package 9 / address 12) is the one that actually has the wrong
# address; its original (mistaken) available_time now belongs in the new
# address_available_time column, and its available_time reverts to the
# normal 8:00 AM start used elsewhere.
$pkg.Range("I10").Value2 = $pkg.Range("C10").Value2
$pkg.Range("C10").Value2 = 0.33333333333333331
$pkg.Range("H10").Value2 = $true

# --- Active sheet / selection bookkeeping -----------------------------------
# Previously "distance_data" was the tab shown on open; now "package_data"
# is shown on open, scrolled back to the top with cell N3 selected.
$dist.Activate()
$dist.Range("A1").Select() | Out-Null

$pkg.Activate()
$pkg.Range("N3").Select() | Out-Null
